$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4
$ws.Range("C3").Value = 4
$ws.Range("C4").Value = 6.5
$ws.Range("C5").Value = 7.5
$ws.Range("C6").Value = 3
$ws.Range("C7").Value = 8

$ws.Range("H22").Value = "old"
$ws.Range("H23").Value = 0
$ws.Range("H24").Value = 0
$ws.Range("H25").Value = 8.67
$ws.Range("H26").Value = 2.27
$ws.Range("H27").Value = 6.53
$ws.Range("H28").Value = 0.2

[void]$ws.Range("D12").Select()
